$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Cells.Item(1,1).Value = 'venue'
$ws.Cells.Item(1,2).Value = 'date'
$ws.Cells.Item(1,3).Value = 'result'
$ws.Cells.Item(1,4).Value = 'ownTeam'
$ws.Cells.Item(1,5).Value = 'oppTeam'
$ws.Cells.Item(1,6).Value = 'batsman'
$ws.Cells.Item(1,7).Value = 'totalRuns'
$ws.Cells.Item(1,8).Value = 'totalBalls'
$ws.Cells.Item(1,9).Value = 'total4s'
$ws.Cells.Item(1,10).Value = 'total6s'
$ws.Cells.Item(1,11).Value = 'sr'

# --- Data rows ---
$ws.Cells.Item(2,1).Value = ' Dubai (DSC)'
$ws.Cells.Item(2,2).Value = ' September 21 2020'
$ws.Cells.Item(2,3).Value = 'RCB won by 10 runs'
$ws.Cells.Item(2,4).Value = 'Sunrisers Hyderabad'
$ws.Cells.Item(2,5).Value = 'Royal Challengers Bangalore'
$ws.Cells.Item(2,6).Value = 'Abhishek Sharma '
$ws.Cells.Item(2,7).Formula = '="7"'
$ws.Cells.Item(2,8).Formula = '="4"'
$ws.Cells.Item(2,9).Formula = '="1"'
$ws.Cells.Item(2,10).Formula = '="0"'
$ws.Cells.Item(2,11).Formula = '="175.00"'
$ws.Cells.Item(3,1).Value = ' Sharjah'
$ws.Cells.Item(3,2).Value = ' October 31 2020'
$ws.Cells.Item(3,3).Value = 'Sunrisers won by 5 wickets (with 35 balls remaining)'
$ws.Cells.Item(3,4).Value = 'Sunrisers Hyderabad'
$ws.Cells.Item(3,5).Value = 'Royal Challengers Bangalore'
$ws.Cells.Item(3,6).Value = 'Abhishek Sharma '
$ws.Cells.Item(3,7).Formula = '="8"'
$ws.Cells.Item(3,8).Formula = '="5"'
$ws.Cells.Item(3,9).Formula = '="0"'
$ws.Cells.Item(3,10).Formula = '="1"'
$ws.Cells.Item(3,11).Formula = '="160.00"'
$ws.Cells.Item(4,1).Value = ' Abu Dhabi'
$ws.Cells.Item(4,2).Value = ' September 26 2020'
$ws.Cells.Item(4,3).Value = 'KKR won by 7 wickets (with 12 balls remaining)'
$ws.Cells.Item(4,4).Value = 'Sunrisers Hyderabad'
$ws.Cells.Item(4,5).Value = 'Kolkata Knight Riders'
$ws.Cells.Item(4,6).Value = 'Abhishek Sharma '
$ws.Cells.Item(4,7).Formula = '="2"'
$ws.Cells.Item(4,8).Formula = '="3"'
$ws.Cells.Item(4,9).Formula = '="0"'
$ws.Cells.Item(4,10).Formula = '="0"'
$ws.Cells.Item(4,11).Formula = '="66.66"'
$ws.Cells.Item(5,1).Value = ' Dubai (DSC)'
$ws.Cells.Item(5,2).Value = ' October 08 2020'
$ws.Cells.Item(5,3).Value = 'Sunrisers won by 69 runs'
$ws.Cells.Item(5,4).Value = 'Sunrisers Hyderabad'
$ws.Cells.Item(5,5).Value = 'Kings XI Punjab'
$ws.Cells.Item(5,6).Value = 'Abhishek Sharma '
$ws.Cells.Item(5,7).Formula = '="12"'
$ws.Cells.Item(5,8).Formula = '="6"'
$ws.Cells.Item(5,9).Formula = '="1"'
$ws.Cells.Item(5,10).Formula = '="1"'
$ws.Cells.Item(5,11).Formula = '="200.00"'
$ws.Cells.Item(6,1).Value = ' Dubai (DSC)'
$ws.Cells.Item(6,2).Value = ' October 02 2020'
$ws.Cells.Item(6,3).Value = 'Sunrisers won by 7 runs'
$ws.Cells.Item(6,4).Value = 'Sunrisers Hyderabad'
$ws.Cells.Item(6,5).Value = 'Chennai Super Kings'
$ws.Cells.Item(6,6).Value = 'Abhishek Sharma '
$ws.Cells.Item(6,7).Formula = '="31"'
$ws.Cells.Item(6,8).Formula = '="24"'
$ws.Cells.Item(6,9).Formula = '="4"'
$ws.Cells.Item(6,10).Formula = '="1"'
$ws.Cells.Item(6,11).Formula = '="129.16"'
$ws.Cells.Item(7,1).Value = ' Sharjah'
$ws.Cells.Item(7,2).Value = ' October 04 2020'
$ws.Cells.Item(7,3).Value = 'Mumbai won by 34 runs'
$ws.Cells.Item(7,4).Value = 'Sunrisers Hyderabad'
$ws.Cells.Item(7,5).Value = 'Mumbai Indians'
$ws.Cells.Item(7,6).Value = 'Abhishek Sharma '
$ws.Cells.Item(7,7).Formula = '="10"'
$ws.Cells.Item(7,8).Formula = '="13"'
$ws.Cells.Item(7,9).Formula = '="0"'
$ws.Cells.Item(7,10).Formula = '="0"'
$ws.Cells.Item(7,11).Formula = '="76.92"'
$ws.Cells.Item(8,1).Value = ' Abu Dhabi'
$ws.Cells.Item(8,2).Value = ' September 29 2020'
$ws.Cells.Item(8,3).Value = 'Sunrisers won by 15 runs'
$ws.Cells.Item(8,4).Value = 'Sunrisers Hyderabad'
$ws.Cells.Item(8,5).Value = 'Delhi Capitals'
$ws.Cells.Item(8,6).Value = 'Abhishek Sharma '
$ws.Cells.Item(8,7).Formula = '="1"'
$ws.Cells.Item(8,8).Formula = '="1"'
$ws.Cells.Item(8,9).Formula = '="0"'
$ws.Cells.Item(8,10).Formula = '="0"'
$ws.Cells.Item(8,11).Formula = '="100.00"'

# --- Convert numeric-looking text formulas to static text values (avoid Excel auto-typing) ---
$ws.Range("A1:K8").Copy()
$ws.Range("A1:K8").PasteSpecial(-4163)

Write-Host "Done"